$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-29 Tuesday" "2024-10-30 Wednesday"

Replace-Text "498÷5=" "183÷9="
Replace-Text "649÷4=" "147÷4="
Replace-Text "906÷8=" "530÷9="
Replace-Text "235÷2=" "591÷4="
Replace-Text "591÷9=" "287÷5="
Replace-Text "172÷5=" "719÷2="
Replace-Text "990÷2=" "922÷5="
Replace-Text "388÷6=" "274÷2="
Replace-Text "802÷8=" "645÷2="
Replace-Text "614÷3=" "226÷8="
Replace-Text "924÷4=" "284÷7="
Replace-Text "420÷5=" "911÷4="
Replace-Text "699÷3=" "744÷5="
Replace-Text "156÷7=" "903÷4="
Replace-Text "386÷7=" "342÷9="
Replace-Text "569÷5=" "114÷6="
Replace-Text "371÷6=" "666÷4="
Replace-Text "358÷2=" "367÷8="
Replace-Text "909÷7=" "970÷7="
Replace-Text "844÷9=" "854÷9="
Replace-Text "605÷7=" "920÷3="
Replace-Text "850÷4=" "995÷5="
Replace-Text "439÷3=" "216÷3="
Replace-Text "895÷2=" "245÷4="
Replace-Text "124÷4=" "329÷5="
